$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Restructure rows FIRST (before creating any helper cells, so nothing we
#    place gets shifted by the row insert/delete):
#    insert a new row at 5 (old row5 "Source" -> row6, old row6 "Note" -> row7),
#    then delete the old "Note" row entirely.
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(7).Delete()

# ---------------------------------------------------------------------------
# 2. Capture template formats from existing cells by copying their formats
#    into holding cells far away on the sheet. This lets Excel reuse/derive
#    style indices the same way it would when a human duplicates formatting.
#    (row 4 is still the original "header" row at this point, and row 6 is
#    the original "Source" row, both untouched so far)
# ---------------------------------------------------------------------------
# Template for "text label" cells in rows 4/5 (old A4 had font Arial10/theme1,
# fill, top+bottom border, left/center/wrap alignment)
$ws.Range("A4").Copy()
$ws.Range("M4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> template (top border variant)
$ws.Range("A4").Copy()
$ws.Range("M5").PasteSpecial(-4122) | Out-Null   # template (will become bottom border variant)

# Template for numeric data cells in rows 4/5 (old B4 had font Arial10/indexed8,
# fill, numFmt 164, no border, right aligned)
$ws.Range("B4").Copy()
$ws.Range("N4").PasteSpecial(-4122) | Out-Null

# Template for the Source-row label cell (old row6 B:H used font Arial9/indexed8,
# fill, top border, left/center/wrap)
$ws.Range("B6").Copy()
$ws.Range("O4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Row 1 - Title
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Lanchkhuti Municipality"
$ws.Range("A1:I1").Merge()
$ws.Rows.Item(1).RowHeight = 51
$titleRng = $ws.Range("A1:I1")
$titleRng.HorizontalAlignment = -4108  # xlCenter
$titleRng.VerticalAlignment = -4108    # xlCenter
$titleRng.WrapText = $true

# ---------------------------------------------------------------------------
# 4. Row 3 - A3 font becomes Sylfaen 11 (keeps its top border)
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# ---------------------------------------------------------------------------
# 5. Row 4 - "family with disabilities Persons" data row
# ---------------------------------------------------------------------------
$ws.Range("M4").Copy()
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none (top border only remains)
$ws.Rows.Item(4).RowHeight = 24.75

$cols = @("B","C","D","E","F","G","H","I")
$nums4 = @(1293,1258,1188,1230,1216,3227,1119,1099)
for ($i = 0; $i -lt 8; $i++) {
  $addr = "$($cols[$i])4"
  $ws.Range("N4").Copy()
  $ws.Range($addr).PasteSpecial(-4122) | Out-Null
  $ws.Range($addr).Value = $nums4[$i]
  $ws.Range($addr).HorizontalAlignment = 1  # xlGeneral
}

# ---------------------------------------------------------------------------
# 6. Row 5 - "disabilities Persons" data row
# ---------------------------------------------------------------------------
$ws.Range("M5").Copy()
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Borders.Item(8).LineStyle = -4142   # xlEdgeTop -> none
$ws.Range("A5").Borders.Item(9).LineStyle = 1       # xlEdgeBottom -> thin
$ws.Range("A5").Borders.Item(9).Weight = 2
$ws.Rows.Item(5).RowHeight = 21

$nums5 = @(1554,1512,1439,1497,1475,3718,1340,1316)
for ($i = 0; $i -lt 8; $i++) {
  $addr = "$($cols[$i])5"
  $ws.Range("N4").Copy()
  $ws.Range($addr).PasteSpecial(-4122) | Out-Null
  $ws.Range($addr).Value = $nums5[$i]
  $ws.Range($addr).HorizontalAlignment = 1  # xlGeneral
}
# I5 additionally gets a bottom border like A5
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------------
# 7. Row 6 - Source row (shifted down from old row 5). Restyle A6 to use the
#    non-bold 9pt font, add fill + left/center/wrap alignment, no border.
# ---------------------------------------------------------------------------
$ws.Range("O4").Copy()
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("A6").Borders.Item(9).LineStyle = -4142  # xlEdgeBottom none
$ws.Range("A6").Borders.Item(8).LineStyle = -4142  # xlEdgeTop none
$ws.Rows.Item(6).RowHeight = 27.75

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 8. Clean up helper/template cells used for format-copying.
# ---------------------------------------------------------------------------
$ws.Range("M4:O5").Clear()

# ---------------------------------------------------------------------------
# 9. Column A width
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.96

Write-Host "Done"
for ($r = 1; $r -le 6; $r++) {
  $rowVals = @()
  for ($c = 1; $c -le 9; $c++) {
    $rowVals += $ws.Cells.Item($r,$c).Value2
  }
  Write-Host "Row $r : $($rowVals -join ' | ')"
}
Write-Host "UsedRange: $($ws.UsedRange.Address())"
